$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "60.148.89"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "2.630.88"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.34"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +6.97%  "
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("E10").Value = "  +3.34%  "
$ws.Range("E11").Value = "  +6.37%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").Value = "3.095.96"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "60.132.36"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("D16").Value = "2.645.48"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("E18").Value = "  +4.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.436"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.98%  "
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("D28").Value = "0.0₃0771"
$ws.Range("E28").Value = "  +4.21%  "
$ws.Range("E30").Value = "  +4.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.15"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.914"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +10.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.914"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.27%  "
$ws.Range("E37").Value = "  +5.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  +5.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "303.50"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.62%  "
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.604"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0977"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.16%  "
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.70"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.31%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.34"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +9.84%  "
$ws.Range("D51").Value = "1.957.66"
$ws.Range("E51").Value = "  +0.85%  "
